{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// The commit joins three sentence fragments that had been split across\n// multiple <w:r> runs back into single runs (no visible text changes),\n// and removes a Wingdings smiley (and the stray trailing space that\n// followed it) from the closing \"Let me know what you think\" line.\n//\n// Strategy: locate each affected span of text with a Range search, then\n// re-insert the same (or trimmed) text as a \"Replace\" on that range.\n// Replacing a multi-run span with literal text collapses it back down\n// to a single run, which mirrors the run-merging shown in the diff.\n\nconst body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\n// \" and nutritious. \" + \"It is a seasonal favorite for my friends and\n// family\" + \" and I usually the ingredients for under 20 dollars TOTAL\n// (not including the spices and oil)\" -> merged into a single run.\nconst part1 =\n  \" and nutritious. It is a seasonal favorite for my friends and family\" +\n  \" and I usually the ingredients for under 20 dollars TOTAL (not including the spices and oil)\";\nconst found1 = body.search(part1, { matchCase: true });\nfound1.load(\"items\");\nawait context.sync();\nif (found1.items.length > 0) {\n  found1.items[0].insertText(part1, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2 -----------------------------------------------------------\n// \"If you can\\u2019t get your hands on these ingredients or don\\u2019t\n// like some, don\\u2019t sweat it. \" + \"All of these are easily\n// replaceable, and you can make the dish according to your tastes and\n// need\" + \"s. \" -> merged into a single run.\nconst part2 =\n  \"If you can\\u2019t get your hands on these ingredients or don\\u2019t like some, don\\u2019t sweat it. \" +\n  \"All of these are easily replaceable, and you can make the dish according to your tastes and needs. \";\nconst found2 = body.search(part2, { matchCase: true });\nfound2.load(\"items\");\nawait context.sync();\nif (found2.items.length > 0) {\n  found2.items[0].insertText(part2, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 3 -----------------------------------------------------------\n// Remove the Wingdings smiley symbol run (and the lone space run that\n// followed it) at the end of \"Let me know what you think <sym> \".\n// The <w:sym> element carries no character in Range/search text, so\n// target the surrounding \"think  \" (two trailing spaces, one on each\n// side of the symbol run) and collapse it down to \"think \" (one space).\nconst part3 = \"think  \";\nconst found3 = body.search(part3, { matchCase: true });\nfound3.load(\"items\");\nawait context.sync();\nif (found3.items.length > 0) {\n  found3.items[0].insertText(\"think \", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n#\n# The commit joins three sentence fragments that had been split across\n# multiple runs back into single runs (no visible text changes), and\n# removes a Wingdings smiley (plus the stray trailing space that\n# followed it) from the closing \"Let me know what you think\" line.\n#\n# Strategy: use Range.Find.Execute(... Replace:=wdReplaceOne) to replace\n# each affected span of text with itself (or trimmed text). Word\n# collapses a multi-run span that is replaced by literal text back down\n# into a single run, which mirrors the run-merging shown in the diff.\n\n$d = $word.ActiveDocument\n\n# --- Change 1 -------------------------------------------------------------\n# \" and nutritious. \" + \"It is a seasonal favorite for my friends and\n# family\" + \" and I usually the ingredients for under 20 dollars TOTAL\n# (not including the spices and oil)\" -> merged into a single run.\n$part1 = \" and nutritious. It is a seasonal favorite for my friends and family and I usually the ingredients for under 20 dollars TOTAL (not including the spices and oil)\"\n$rng1 = $d.Content\n$rng1.Find.Execute($part1, $false, $false, $false, $false, $false, $true, 1, $false, $part1, 2) | Out-Null\n\n# --- Change 2 -------------------------------------------------------------\n# \"If you can't get your hands on these ingredients or don't like some,\n# don't sweat it. \" + \"All of these are easily replaceable, and you can\n# make the dish according to your tastes and need\" + \"s. \" -> merged\n# into a single run.\n$apos = [char]0x2019\n$part2 = \"If you can\" + $apos + \"t get your hands on these ingredients or don\" + $apos + \"t like some, don\" + $apos + \"t sweat it. All of these are easily replaceable, and you can make the dish according to your tastes and needs. \"\n$rng2 = $d.Content\n$rng2.Find.Execute($part2, $false, $false, $false, $false, $false, $true, 1, $false, $part2, 2) | Out-Null\n\n# --- Change 3 ---------------------------------------------------------------\n# Remove the Wingdings smiley symbol run (and the lone space run that\n# followed it) at the end of \"Let me know what you think <sym> \". The\n# <w:sym> element carries no character in the Range text, so target the\n# surrounding \"think  \" (two trailing spaces, one on each side of the\n# symbol run) and collapse it down to \"think \" (one space), which drops\n# the symbol run and the extra space run along with it.\n$part3 = \"think  \"\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute($part3)\nif ($found3) {\n    $rng3.Text = \"think \"\n}\n"}
